# Update countries & provincias Spain
# Applies the 27-May-2020 15:05 -> 15:35 data refresh to the "Pais" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header timestamp -------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 27 de Mayo de 2020 a las 15:35"

# --- Row 4: Estados Unidos (in-place stat refresh, no reordering) -----
$ws.Range("B4").Value = 1727992
$ws.Range("C4").Value = 2717
$ws.Range("D4").Value = 480133
$ws.Range("E4").Value = 1147230
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 57
$ws.Range("H4").Value = 100629

# --- Rows 100-102: Kenia gets fresh numbers and moves ahead of --------
# --- Eslovenia / Maldivas (which keep their own numbers, shifted) -----
$ws.Range("A100").Value = "Kenia"
$ws.Range("B100").Value = 1471
$ws.Range("C100").Value = 123
$ws.Range("D100").Value = 408
$ws.Range("E100").Value = 1008
$ws.Range("F100").Value = 0
$ws.Range("G100").Value = 3
$ws.Range("H100").Value = 55

$ws.Range("A101").Value = "Eslovenia"
$ws.Range("B101").Value = 1471
$ws.Range("C101").Value = 2
$ws.Range("D101").Value = 1354
$ws.Range("E101").Value = 9
$ws.Range("F101").Value = 0
$ws.Range("G101").Value = 0
$ws.Range("H101").Value = 108

$ws.Range("A102").Value = "Maldivas"
$ws.Range("B102").Value = 1457
$ws.Range("C102").Value = 19
$ws.Range("D102").Value = 197
$ws.Range("E102").Value = 1255
$ws.Range("F102").Value = 0
$ws.Range("G102").Value = 0
$ws.Range("H102").Value = 5

# --- Rows 199-201: Santa Lucia moves ahead of Belice / Nueva --------
# --- Caledonia; all three keep their own numbers, just shifted ------
$ws.Range("A199").Value = "Santa Lucia"
$ws.Range("B199").Value = 18
$ws.Range("C199").Value = 0
$ws.Range("D199").Value = 18
$ws.Range("E199").Value = 0
$ws.Range("F199").Value = 0
$ws.Range("G199").Value = 0
$ws.Range("H199").Value = 0

$ws.Range("A200").Value = "Belice"
$ws.Range("B200").Value = 18
$ws.Range("C200").Value = 0
$ws.Range("D200").Value = 16
$ws.Range("E200").Value = 0
$ws.Range("F200").Value = 0
$ws.Range("G200").Value = 0
$ws.Range("H200").Value = 2

$ws.Range("A201").Value = "Nueva Caledonia"
$ws.Range("B201").Value = 18
$ws.Range("C201").Value = 0
$ws.Range("D201").Value = 18
$ws.Range("E201").Value = 0
$ws.Range("F201").Value = 0
$ws.Range("G201").Value = 0
$ws.Range("H201").Value = 0

# --- Rows 207-208: Groenlandia swaps ahead of Islas Turcas y Caicos --
$ws.Range("A207").Value = "Groenlandia"
$ws.Range("B207").Value = 12
$ws.Range("C207").Value = 0
$ws.Range("D207").Value = 11
$ws.Range("E207").Value = 1
$ws.Range("F207").Value = 0
$ws.Range("G207").Value = 0
$ws.Range("H207").Value = 0

$ws.Range("A208").Value = "Islas Turcas y Caicos"
$ws.Range("B208").Value = 12
$ws.Range("C208").Value = 0
$ws.Range("D208").Value = 10
$ws.Range("E208").Value = 1
$ws.Range("F208").Value = 0
$ws.Range("G208").Value = 0
$ws.Range("H208").Value = 1

# --- Rows 210-211: Montserrat swaps ahead of Seychelles --------------
$ws.Range("A210").Value = "Montserrat"
$ws.Range("B210").Value = 11
$ws.Range("C210").Value = 0
$ws.Range("D210").Value = 10
$ws.Range("E210").Value = 0
$ws.Range("F210").Value = 0
$ws.Range("G210").Value = 0
$ws.Range("H210").Value = 1

$ws.Range("A211").Value = "Seychelles"
$ws.Range("B211").Value = 11
$ws.Range("C211").Value = 0
$ws.Range("D211").Value = 11
$ws.Range("E211").Value = 0
$ws.Range("F211").Value = 0
$ws.Range("G211").Value = 0
$ws.Range("H211").Value = 0
